$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010796785354614
$ws.Range("B1").Value = 2.124924182891846
$ws.Range("C1").Value = 6.002388954162598
$ws.Range("D1").Value = 1.225194931030273
$ws.Range("E1").Value = 1.230377316474915
